$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A to fit the new longer label.
# NOTE: ColumnWidth is quantized by this runtime to an integer pixel grid
# (1/6-character steps), so an input of 22.0 is the closest achievable value
# that reproduces the target stored width of 22.83203125 (-> 22.8333...).
$ws.Columns.Item(1).ColumnWidth = 22

# Add "MONTEREY AREA TOTALS" label into A2 (matching the style used by the
# other port-name cells in column A, e.g. A4: Calibri 12) and change B2 from
# "MONTEREY AREA TOTALS" to "Totals"
$ws.Range("A2").Value = "MONTEREY AREA TOTALS"
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 12

$ws.Range("B2").Value = "Totals"

# Update the active selection to B7
$ws.Range("B7").Select()
